$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Input": record the source data filename as extra context for this
# analysis, in a newly widened column D.
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Columns.Item(4).ColumnWidth = 66.7109375
$wsInput.Range("D3").Value = "Filename: C:\Neptune\User\Neptune\Data\UTh\2015\0815\011_7184.dat"

# ---------------------------------------------------------------------------
# Sheet "Calc": repeated analyses are now saved separately, so the Monte
# Carlo success-rate columns are appended, and the refined age/error values
# (recomputed from the repeated run) are updated.
# ---------------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("Calc")

# New columns BG / BH with headers + units
$wsCalc.Range("BG1").Value = "Unkorr. Montefehler Erfolgsrate"
$wsCalc.Range("BG2").Value = "(%)"
$wsCalc.Range("BG3").Value = 100

$wsCalc.Range("BH1").Value = "Korr. Montefehler Erfolgsrate"
$wsCalc.Range("BH2").Value = "(%)"
$wsCalc.Range("BH3").Value = 100

# Column widths
$wsCalc.Columns.Item(49).ColumnWidth = 9.7109375
$wsCalc.Columns.Item(50).ColumnWidth = 19.7109375
$wsCalc.Columns.Item(51).ColumnWidth = 19.7109375
$wsCalc.Columns.Item(55).ColumnWidth = 19.7109375
$wsCalc.Columns.Item(59).ColumnWidth = 32.7109375
$wsCalc.Columns.Item(60).ColumnWidth = 30.7109375

# Updated values (recalculated results)
$wsCalc.Range("AP3").Value = 0.534
$wsCalc.Range("AQ3").Value = 0.1890292921065403
$wsCalc.Range("AW3").Value = 0.5303
$wsCalc.Range("AX3").Value = 0.5347190807924826
$wsCalc.Range("AY3").Value = 0.1877546309327436
$wsCalc.Range("BC3").Value = 0.5689992644482499
$wsCalc.Range("BE3").Value = 267.3595403962413
$wsCalc.Range("BF3").Value = 0.189319222453118

# ---------------------------------------------------------------------------
# Sheet "Results": mirrors the updated age-error values from Calc, and a
# couple of column widths change to better fit the refreshed content.
# ---------------------------------------------------------------------------
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Columns.Item(16).ColumnWidth = 8.7109375
$wsResults.Columns.Item(18).ColumnWidth = 19.7109375

$wsResults.Range("N3").Value = 0.534
$wsResults.Range("P3").Value = 0.5303
$wsResults.Range("R3").Value = 0.5689992644482499

# ---------------------------------------------------------------------------
# Sheet "Constants": the 230/229 ratio slope constant is corrected.
# ---------------------------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("B3").Value = 0.00005
